# code fix exam and hsr
# Updates sample/test data rows on the V2Project, Exam and Hsr sheets.

$wb = $excel.ActiveWorkbook

# ---- V2Project sheet: CardHolderName (P2) ----
$wsV2 = $wb.Worksheets.Item("V2Project")
$wsV2.Range("P2").Value = "Norris"

# ---- Exam sheet: refresh row 2 sample data ----
$wsExam = $wb.Worksheets.Item("Exam")
# Force text format first for purely-numeric-looking values so the
# leading-digit strings are preserved exactly (not auto-converted to numbers).
$wsExam.Range("B2").NumberFormat = "@"
$wsExam.Range("B2").Value = "1830104860"
$wsExam.Range("C2").Value = "United States"
$wsExam.Range("D2").Value = "Montana"
$wsExam.Range("E2").Value = "6935 Lizette Ridges"
$wsExam.Range("F2").Value = "West Marinda"
$wsExam.Range("G2").NumberFormat = "@"
$wsExam.Range("G2").Value = "95210"
$wsExam.Range("H2").Value = "Wilburn"
$wsExam.Range("I2").Value = "100 RESILIENT CITIES"

# ---- Hsr sheet: refresh row 2 sample data ----
$wsHsr = $wb.Worksheets.Item("Hsr")
$wsHsr.Range("A2").Value = "HSRP006238"
$wsHsr.Range("B2").Value = "Automation HSR Project 7162135"
$wsHsr.Range("C2").Value = "1-800 Flowers.com Inc."
$wsHsr.Range("D2").Value = "Manufacturing"
$wsHsr.Range("F2").Value = "Kansas"
$wsHsr.Range("G2").Value = "074 Krystyna Shoal"
$wsHsr.Range("H2").Value = "Chantayburgh"
$wsHsr.Range("I2").Value = "00030-5564"
# K2 has a significant leading zero ("015332"); force text so it is not
# silently turned into the number 15332.
$wsHsr.Range("K2").NumberFormat = "@"
$wsHsr.Range("K2").Value = "015332"
